$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AZ holds "recentproj_yr" - the most recent year data was projected/updated.
# Update it from 2015 to 2016 for the rows whose recent-update year has moved forward.
$rowsToUpdate = @(2, 4, 5, 6, 7, 8, 10, 11, 13, 16, 17, 19, 20)

foreach ($r in $rowsToUpdate) {
    $ws.Range("AZ$r").Value = 2016
}
